# "Add wool to price variation" -- update the existing scenario rows 2-5
# with refreshed simulated values, and append four new scenario rows
# (6-9, labelled c1_4 .. c1_7) to each of the grain / meat / wool / prob
# sheets.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "grain" (sheet1): columns B:K repeat the same value across a row.
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("grain")

$ws.Range("B2:K2").Value = 0.8917796022511414
$ws.Range("B3:K3").Value = 0.9097077561761636
$ws.Range("B4:K4").Value = 0.862806697296832
$ws.Range("B5:K5").Value = 0.885528406419913

$ws.Range("A5:K5").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:K7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:K8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:K9").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "c1_4"
$ws.Range("B6:K6").Value = 1.114471593580087

$ws.Range("A7").Value = "c1_5"
$ws.Range("B7:K7").Value = 1.137193302703168

$ws.Range("A8").Value = "c1_6"
$ws.Range("B8:K8").Value = 1.090292243823836

$ws.Range("A9").Value = "c1_7"
$ws.Range("B9:K9").Value = 1.108220397748859

# ----------------------------------------------------------------------
# Sheet "meat" (sheet2)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("meat")

$ws.Range("B2:K2").Value = 0.8040309623267345
$ws.Range("B3:K3").Value = 0.8301305515595439
$ws.Range("B4:K4").Value = 1.215804309908924
$ws.Range("B5:K5").Value = 1.248158231333236

$ws.Range("A5:K5").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:K7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:K8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:K9").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "c1_4"
$ws.Range("B6:K6").Value = 0.7518417686667644

$ws.Range("A7").Value = "c1_5"
$ws.Range("B7:K7").Value = 0.784195690091076

$ws.Range("A8").Value = "c1_6"
$ws.Range("B8:K8").Value = 1.169869448440456

$ws.Range("A9").Value = "c1_7"
$ws.Range("B9:K9").Value = 1.195969037673266

# ----------------------------------------------------------------------
# Sheet "wool" (sheet3)
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("wool")

$ws.Range("B2:K2").Value = 0.9093036898722159
$ws.Range("B3:K3").Value = 1.065598434559732
$ws.Range("B4:K4").Value = 0.921135610673716
$ws.Range("B5:K5").Value = 1.075750587078096

$ws.Range("A5:K5").Copy() | Out-Null
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:K7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:K8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:K9").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "c1_4"
$ws.Range("B6:K6").Value = 0.9242494129219041

$ws.Range("A7").Value = "c1_5"
$ws.Range("B7:K7").Value = 1.078864389326284

$ws.Range("A8").Value = "c1_6"
$ws.Range("B8:K8").Value = 0.9344015654402679

$ws.Range("A9").Value = "c1_7"
$ws.Range("B9:K9").Value = 1.090696310127785

# ----------------------------------------------------------------------
# Sheet "prob" (sheet4) -- only columns A:B
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("prob")

$ws.Range("B2").Value = 0.1300397501860382
$ws.Range("B3").Value = 0.07264965140027922
$ws.Range("B4").Value = 0.1539326632094
$ws.Range("B5").Value = 0.1433779352042823

$ws.Range("A5:B5").Copy() | Out-Null
$ws.Range("A6:B6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7:B7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8:B8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9:B9").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "c1_4"
$ws.Range("B6").Value = 0.1433779352042824

$ws.Range("A7").Value = "c1_5"
$ws.Range("B7").Value = 0.1539326632094002

$ws.Range("A8").Value = "c1_6"
$ws.Range("B8").Value = 0.07264965140027926

$ws.Range("A9").Value = "c1_7"
$ws.Range("B9").Value = 0.1300397501860383
